$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.423.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3819"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.222"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.417"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.307"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001234"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.628.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06960"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.578"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.404.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.536"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.069"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.268"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.810.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.090"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.154"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.506"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02757"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2501"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08752"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07022"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("E41").Value = "  -1.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6462"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.966"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07929"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
